$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timesheet")

# Activity description on row 9 changed
$ws.Range("B9").Value = "Development ( 1 ) "

# Daily hours for the activity row (row 9) -- all become the text "1".
# Using TEXT()+PasteSpecial(values) keeps the cells as shared-string text
# (t="s") with their original style/number-format, instead of letting a
# literal numeric-looking string get auto-coerced to a numeric cell.
$ws.Range("E9:K9").Formula = "=TEXT(1,""0"")"
$ws.Range("E9:K9").Copy()
$ws.Range("E9:K9").PasteSpecial(-4163)

# Mirrored totals row (row 18) -- all become "1" as well
$ws.Range("E18:K18").Formula = "=TEXT(1,""0"")"
$ws.Range("E18:K18").Copy()
$ws.Range("E18:K18").PasteSpecial(-4163)

# Total hours for the week
$ws.Range("D19").Formula = "=TEXT(7,""0"")"
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)

$excel.CutCopyMode = 0
